# Applies the "Synthèse projet MNT" revision described by the commit diff.
# Runs against $word.ActiveDocument (iron_native headless Word COM-interop).

$d = $word.ActiveDocument

function ReplaceOnce($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND:" $find
    }
    return $ok
}

# 1) "...apparu complexe, en effet les méthodes..." -> split into two
#    sentences and extend the paragraph with a new closing clause.
ReplaceOnce `
    "complexe, en effet les méthodes à mettre en place n’étaient pas définies explicitement." `
    "complexe. En effet les méthodes à mettre en place n’étaient pas définies explicitement, il nous a fallu faire beaucoup d’hypothèses quant aux attentes vis-à-vis de ce projet. "

# 2) "Dans le but de nous organiser correctement au fur et à mesure du projet
#    nous avions décidé..." -> "...organiser efficacement nous avions décidé..."
ReplaceOnce `
    "Dans le but de nous organiser correctement au fur et à mesure du projet nous avions décidé" `
    "Dans le but de nous organiser efficacement nous avions décidé"

# 3) "De plus nous avons, autant que possible, ... mais aussi cela nous
#    assurait de ne rien oublier ..." -> "Nous avons, ... mais cela nous
#    assurait aussi de ne rien oublier ..."
ReplaceOnce `
    "De plus nous avons, autant que possible, rédigé les rapports en parallèle de la programmation. De cette façon nous étions sûr de ne pas prendre de retard dans la rédaction de ces rapports, mais aussi cela nous assurait de ne rien oublier d’important pour le bilan du projet." `
    "Nous avons, autant que possible, rédigé les rapports en parallèle de la programmation. De cette façon nous étions sûr de ne pas prendre de retard dans la rédaction de ces rapports, mais cela nous assurait aussi de ne rien oublier d’important pour le bilan du projet."

# 4) "... sont toujours présents même s'ils ne représentent ..." -> add
#    "à titre informatif" before "même s'ils".
ReplaceOnce `
    "sont toujours présents même s’ils ne représentent plus forcément d’utilité." `
    "sont toujours présents à titre informatif même s’ils ne représentent plus forcément d’utilité."

# 5) Heading "II. Réussites et échecs" -> "II. Réussites et améliorations"
ReplaceOnce "II. Réussites et échecs" "II. Réussites et améliorations"

# 7) "...un nuage de points de 3 coordonnées (XYZ). Les points ne sont pas
#    forcément séparés..." -> "... (XYZ) et que les points ne sont pas
#    forcément séparés..."
ReplaceOnce `
    "un nuage de points de 3 coordonnées (XYZ). Les points ne sont pas forcément séparés" `
    "un nuage de points de 3 coordonnées (XYZ) et que les points ne sont pas forcément séparés"

# 6) The paragraph starting "Nous sommes partie avec l'hypothèse..." gains a
#    first-line indent (0.5cm / 708 twips == 35.4pt).
$rng = $d.Content
$rng.Find.Execute("Nous sommes partie avec l’hypothèse")
if ($rng.Find.Found) {
    $rng.ParagraphFormat.FirstLineIndent = 35.4
} else {
    Write-Host "NOT FOUND: Nous sommes partie avec l'hypothese paragraph"
}

# 8) The paragraph starting "Nous devions pouvoir visualiser..." also gains
#    the same first-line indent.
$rng = $d.Content
$rng.Find.Execute("Nous devions pouvoir visualiser le MNT")
if ($rng.Find.Found) {
    $rng.ParagraphFormat.FirstLineIndent = 35.4
} else {
    Write-Host "NOT FOUND: Nous devions pouvoir visualiser paragraph"
}

Write-Host "Simple replacements done"
